# BudgetEstimateGridGeomInC++_phase3.xlsx edit:
#  - remove the now-unused "api calls" and "Sheet3" worksheets
#  - the "Total hours spent so far" summary row moves from row 36 to row 38
#    (two extra blank rows were inserted above it)
#  - update the remembered selection on the "begroting" sheet

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("begroting")

# Relocate the trailing "Total hours spent so far" summary (was row 36) down
# to row 38, preserving its formulas (relative refs re-target row 38).
$ws.Range("A38").Value = $ws.Range("A36").Value()
$ws.Range("B38").Formula = "=SUM(B27:B34)"
$ws.Range("C38").Formula = "=B38*135"
$ws.Range("A36:C36").ClearContents()

# Drop the helper sheets that are no longer needed; this also prunes the
# shared-string entries that only they referenced ("FORTRAN subroutines",
# "gridtonet", "API calls", "mergenodes").
$wb.Worksheets.Item("api calls").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Match the saved selection state.
$ws.Range("C30").Select()
